# Atualização automática: 2025-08-31 21:00:26
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 18 (corrected detection image / bbox / confidence) ---
$ws.Range("D18").Value = "image_20250808100711_ppp0.jpg"
$ws.Range("I18").Value = "1182,409,1232,451"
$ws.Range("J18").Value = "'0.75"
$ws.Range("J18").Style = "Normal"

# --- Append new row 23 (new detection record) ---
$ws.Range("A23").Value = "00deb925-04cf-4a0c-b2a2-5289a118de4d"
$ws.Range("B23").Value = "mosca"
$ws.Range("C23").Value = 45900
$ws.Range("C23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D23").Value = "image_20250831214359_ppp0.jpg"
$ws.Range("E23").Value = "PLACA_20250717165933"
$ws.Range("F23").Value = "Beja"
$ws.Range("G23").Value = 38.02035
$ws.Range("H23").Value = -7.94715
$ws.Range("I23").Value = "1256,526,1294,578"
$ws.Range("J23").Value = "'0.72"
$ws.Range("J23").Style = "Normal"
